$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new candidate's data (keep F2 and H2 unchanged)
$ws.Range("B2").Value = 'Kelvin Ee'
$ws.Range("C2").Value = '+60 11-3919 0131'
$ws.Range("D2").Value = 'kelvinee2001@gmail.com'
$ws.Range("E2").Value = 'N/A'
$ws.Range("G2").Value = '[{''job_title'': ''PHP Web Developer'', ''job_company'': ''Powerec Services Sdn Bhd'', ''Industries'': ''N/A'', ''start_date'': ''2021-11'', ''end_date'': ''2022-02'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}, {''job_title'': ''PPG Intern'', ''job_company'': ''N/A'', ''Industries'': ''N/A'', ''start_date'': ''2023-08'', ''end_date'': ''2024-05-19 11:33:54.994675'', ''job_location'': ''N/A'', ''job_duration'': ''N/A''}]'
$ws.Range("I2").Value = '[{''field_of_study'': ''Bachelor of Computer Science (Data Engineering)'', ''level'': "Bachelor''s Degree", ''cgpa'': ''3.98'', ''university'': ''Universiti Teknologi Malaysia'', ''start_date'': ''2020'', ''year_of_graduation'': ''N/A''}]'
$ws.Range("J2").Value = '[''Alteryx Designer Core Certified'']'
$ws.Range("K2").Value = '[''Python'', ''SQL'', ''C++'', ''R'', ''Power BI'', ''Tableau'', ''Alteryx'', ''Excel'', ''Databricks'', ''Azure Data Factory'', ''Blob Storage'', ''KeyVault'', ''HTML'', ''CSS'', ''PHP'', ''Bootstrap'', ''JavaScript'', ''C#'', ''Cloud Foundation'', ''Machine Learning'', ''Data Analytic'']'
$ws.Range("L2").Value = '[''English'', ''Malay'', ''Chinese'']'

# Remove row 3 entirely (the second candidate entry no longer exists)
$ws.Rows("3").Delete()
